# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from serial date 45441 (2024-05-29) to 45442 (2024-05-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45441) {
        $cell.Value2 = 45442
    }
}
